# Economic Dashboard V1 - Update dashboards 2025-12-05
# Applies the cell value / style changes captured in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Style-only changes: move the "Latest Period" highlight (yellow fill,
#    style index 48) on/off specific date cells. We copy formats from a
#    donor cell that already carries the desired style so the workbook's
#    existing style table (cellXfs index 47/48) is reused instead of having
#    new duplicate styles synthesized.
# ---------------------------------------------------------------------------

# Style 48 (yellow "updated" highlight) donor -> apply to C28:C31
$ws.Range("N29").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").PasteSpecial(-4122)

# Style 47 (no highlight) donor -> apply to N39
$ws.Range("N7").Copy()
$ws.Range("N39").PasteSpecial(-4122)

# N51 switches from style 47 -> style 48 (highlighted)
$ws.Range("N29").Copy()
$ws.Range("N51").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row 7 - GDP Nowcast present value refresh
# ---------------------------------------------------------------------------
$ws.Range("F7").Value = 0.3129565816678153

# ---------------------------------------------------------------------------
# 3. Row 28 - Durable Orders M/M (values refreshed; date highlight set above)
# ---------------------------------------------------------------------------
$ws.Range("F28").Value = 0.005052252529329948
$ws.Range("G28").Value = 0.03004963172206243

# ---------------------------------------------------------------------------
# 4. Row 29 - Durable Orders Y/Y + 5yr,5yr Forward inflation series roll
# ---------------------------------------------------------------------------
$ws.Range("F29").Value = 0.07264359641534658
$ws.Range("G29").Value = 0.07661265288383932
$ws.Range("N29").Value = 45995
$ws.Range("R29").Value = 2.19
$ws.Range("S29").Value = 2.18
$ws.Range("T29").Value = 2.17

# ---------------------------------------------------------------------------
# 5. Row 30 - Dur Orders Non Def x Aircraft + 10yr TIPS series roll
# ---------------------------------------------------------------------------
$ws.Range("F30").Value = 0.000905510184407321
$ws.Range("G30").Value = 0.01907672443132968
$ws.Range("N30").Value = 45995
$ws.Range("Q30").Value = 2.26
$ws.Range("T30").Value = 2.24

# ---------------------------------------------------------------------------
# 6. Row 31 - Dur Orders Non Def x Aircraft Y/Y
# ---------------------------------------------------------------------------
$ws.Range("F31").Value = 0.0647360016408333
$ws.Range("G31").Value = 0.06671073894520346

# ---------------------------------------------------------------------------
# 7. Row 47 - FFR latest period date
# ---------------------------------------------------------------------------
$ws.Range("N47").Value = 45994

# ---------------------------------------------------------------------------
# 8. Row 48 - 2y UST series roll
# ---------------------------------------------------------------------------
$ws.Range("N48").Value = 45994
$ws.Range("Q48").Value = 3.49
$ws.Range("R48").Value = 3.51
$ws.Range("S48").Value = 3.54
$ws.Range("U48").ClearContents()

# ---------------------------------------------------------------------------
# 9. Row 49 - 5y UST series roll
# ---------------------------------------------------------------------------
$ws.Range("N49").Value = 45994
$ws.Range("Q49").Value = 3.62
$ws.Range("R49").Value = 3.66
$ws.Range("S49").Value = 3.67
$ws.Range("U49").ClearContents()

# ---------------------------------------------------------------------------
# 10. Row 50 - 10y UST series roll
# ---------------------------------------------------------------------------
$ws.Range("N50").Value = 45994
$ws.Range("Q50").Value = 4.06
$ws.Range("S50").Value = 4.09
$ws.Range("U50").ClearContents()

# ---------------------------------------------------------------------------
# 11. Row 51 - 30y Mortgage series roll (date highlight set above)
# ---------------------------------------------------------------------------
$ws.Range("N51").Value = 45992
$ws.Range("Q51").Value = 6.19
$ws.Range("R51").Value = 6.23
$ws.Range("S51").Value = 6.26
$ws.Range("T51").Value = 6.24
$ws.Range("U51").Value = 6.22

# ---------------------------------------------------------------------------
# 12. Row 52 - BAA series roll
# ---------------------------------------------------------------------------
$ws.Range("N52").Value = 45994
$ws.Range("Q52").Value = 5.83
$ws.Range("R52").Value = 5.85
$ws.Range("S52").Value = 5.87
$ws.Range("U52").ClearContents()
